$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# Document statuses changed: "Vision", "2 fully-dressed use cases" and
# "2 complete use cases with activity diagram" are now complete too, so give
# C3, C5 and C6 the same "complete" checkmark look as C4 (which was already
# marked complete) and then write the checkmark value itself.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C6").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "ü"
$ws.Range("C5").Value = "ü"
$ws.Range("C6").Value = "ü"

# The edit leaves the selection on C6.
$ws.Range("C6").Select() | Out-Null
